$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "Kiefer (n. entb.)"
$ws.Range("C1").Value = "Kiefer (entb.)"
$ws.Range("J1").Value = "Kiefer (n. entb.).letter"
$ws.Range("K1").Value = "Kiefer (entb.).letter"
$ws.Range("N14").Value = "PAR [%]"
